$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper functions: read/write a contiguous block of columns (F..V, i.e. 6..22)
# for a given row, cell by cell (safer than whole-range Value gets/sets in
# this COM host, which flattens multi-cell ranges to a single string).
# NOTE: this engine only supports *positional* function parameters.
# ---------------------------------------------------------------------------
function Get-RowBlock($row, $colStart, $colEnd) {
    $data = @()
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $data += $ws.Cells.Item($row, $c).Value()
    }
    return $data
}

function Set-RowBlock($row, $colStart, $data) {
    $c = $colStart
    foreach ($v in $data) {
        $ws.Cells.Item($row, $c).Value = $v
        $c++
    }
}

# ---------------------------------------------------------------------------
# The match results in columns F:V (home team .. url) were re-shuffled
# between several rows (their A/B/C/D/E "identity" columns - index, country,
# tournament, season, date - stayed put). Apply the same cyclic rotations
# here by reading the old F:V payloads first and then writing them back into
# their new homes.
# ---------------------------------------------------------------------------

# Rows 17 -> 18 -> 19 -> 21 -> 17
$b17 = Get-RowBlock 17 6 22
$b18 = Get-RowBlock 18 6 22
$b19 = Get-RowBlock 19 6 22
$b21 = Get-RowBlock 21 6 22

Set-RowBlock 18 6 $b17
Set-RowBlock 19 6 $b18
Set-RowBlock 21 6 $b19
Set-RowBlock 17 6 $b21

# Rows 85 -> 86 -> 87 -> 88 -> 85
$b85 = Get-RowBlock 85 6 22
$b86 = Get-RowBlock 86 6 22
$b87 = Get-RowBlock 87 6 22
$b88 = Get-RowBlock 88 6 22

Set-RowBlock 86 6 $b85
Set-RowBlock 87 6 $b86
Set-RowBlock 88 6 $b87
Set-RowBlock 85 6 $b88

# Rows 95 -> 97 -> 98 -> 95
$b95 = Get-RowBlock 95 6 22
$b97 = Get-RowBlock 97 6 22
$b98 = Get-RowBlock 98 6 22

Set-RowBlock 97 6 $b95
Set-RowBlock 98 6 $b97
Set-RowBlock 95 6 $b98

# ---------------------------------------------------------------------------
# Append four brand-new match rows (101-104) at the bottom of the sheet.
# Parameters (all positional):
#  1 Row              12 DrawOpenOdds     ...
#  2 Indice            13 DrawOpenData
#  3 DataPartida       14 DrawCloseOdds
#  4 Home              15 DrawCloseData
#  5 HomeGols          16 AwayOpenOdds
#  6 Away              17 AwayOpenData
#  7 AwayGols          18 AwayCloseOdds
#  8 HomeOpenOdds      19 AwayCloseData
#  9 HomeOpenData      20 Url
# 10 HomeCloseOdds
# 11 HomeCloseData
# ---------------------------------------------------------------------------
function Add-MatchRow(
    $Row, $Indice, $DataPartida, $Home, $HomeGols, $Away, $AwayGols,
    $HomeOpenOdds, $HomeOpenData, $HomeCloseOdds, $HomeCloseData,
    $DrawOpenOdds, $DrawOpenData, $DrawCloseOdds, $DrawCloseData,
    $AwayOpenOdds, $AwayOpenData, $AwayCloseOdds, $AwayCloseData, $Url
) {
    $ws.Cells.Item($Row, 1).Value = $Indice
    $ws.Cells.Item($Row, 2).Value = "portugal"
    $ws.Cells.Item($Row, 3).Value = "liga-3"
    $ws.Cells.Item($Row, 4).Value = "2023-2024"
    $ws.Cells.Item($Row, 5).Value = $DataPartida
    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGols
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGols
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenData
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseData
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenData
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseData
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenData
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseData
    $ws.Cells.Item($Row, 22).Value = $Url

    # Match the existing look & feel: column A uses the bold/bordered/centered
    # "index" style, column E uses the date-time number format. Reuse the
    # formatting already present on row 100 rather than constructing new
    # style entries from scratch.
    $ws.Range("A100").Copy() | Out-Null
    $ws.Range("A$Row").PasteSpecial(-4122)
    $ws.Range("E100").Copy() | Out-Null
    $ws.Range("E$Row").PasteSpecial(-4122)
}

Add-MatchRow 101 100 45241.66666666666 `
    "Pero Pinheiro" 1 "Atletico CP" 1 `
    4.91 "07/11/2023 07:11" 5.61 "11/11/2023 15:32" `
    3.84 "07/11/2023 07:11" 3.98 "11/11/2023 15:32" `
    1.71 "07/11/2023 07:11" 1.61 "11/11/2023 15:32" `
    "https://www.betexplorer.com/football/portugal/liga-3/pero-pinheiro-atletico-cp/6axFrHAI/"

Add-MatchRow 102 101 45241.66666666666 `
    "Oliveira Hospital" 0 "Amora" 4 `
    2.03 "05/11/2023 15:44" 2.06 "11/11/2023 15:37" `
    3.32 "05/11/2023 15:44" 3.22 "11/11/2023 15:37" `
    3.81 "05/11/2023 15:44" 4 "11/11/2023 15:37" `
    "https://www.betexplorer.com/football/portugal/liga-3/oliveira-hospital-amora/CxwBqceC/"

Add-MatchRow 103 102 45241.75 `
    "Varzim" 2 "SC Vianense" 0 `
    1.36 "10/11/2023 14:12" 1.76 "11/11/2023 17:56" `
    4.67 "10/11/2023 14:12" 3.79 "11/11/2023 17:56" `
    7.4 "10/11/2023 14:12" 4.61 "11/11/2023 15:56" `
    "https://www.betexplorer.com/football/portugal/liga-3/varzim-sc-vianense/C6dkwtGk/"

Add-MatchRow 104 103 45241.77083333334 `
    "Caldas" 2 "1º Dezembro" 1 `
    1.47 "06/11/2023 15:12" 1.51 "11/11/2023 18:19" `
    4.29 "06/11/2023 15:12" 4.28 "11/11/2023 18:19" `
    6.74 "06/11/2023 15:12" 6.5 "11/11/2023 18:19" `
    "https://www.betexplorer.com/football/portugal/liga-3/caldas-sc-1-dezembro/Yov7pwt6/"

$excel.CutCopyMode = 0
